$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.730.25"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.099.91"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.62"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.27"
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.70"
$ws.Range("E12").Value = "  +5.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.412.02"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.03"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.809"
$ws.Range("E15").Value = "  +3.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.52"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.106.70"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.687.98"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.68"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.12"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.87"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  -3.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.62"
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.85"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.41"
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.31"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  +8.14%  "
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.53"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.02"
$ws.Range("E35").Value = "  +7.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0618"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.12"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.66"
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("E42").Value = "  +3.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.528.40"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("E44").Value = "  +6.51%  "
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.81"
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0911"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.297.45"
$ws.Range("E51").Value = "  +0.48%  "
